$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.083020329475403
$ws.Range("B1").Value = 0.8870450854301453
$ws.Range("C1").Value = 2.505175352096558
$ws.Range("D1").Value = 5.400112628936768
$ws.Range("E1").Value = 1.093494176864624
